$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update toner cost table (columns B and C) for rows 2-15
# per commit: "adequacao tabela de custos do toner por periodo de tempo"

$ws.Range("B2").Value = 12
$ws.Range("C2").Value = 0

$ws.Range("B3").Value = 3
$ws.Range("C3").Value = 12

$ws.Range("B4").Value = 32
$ws.Range("C4").Value = 0

$ws.Range("B5").Value = 4
$ws.Range("C5").Value = 20

$ws.Range("B6").Value = 2
$ws.Range("C6").Value = 8

$ws.Range("B7").Value = 2
$ws.Range("C7").Value = 4

$ws.Range("B8").Value = 6
$ws.Range("C8").Value = 0

$ws.Range("B9").Value = 2
$ws.Range("C9").Value = 4

$ws.Range("C10").Value = 4

$ws.Range("B11").Value = 196
$ws.Range("C11").Value = 0

$ws.Range("B12").Value = 28
$ws.Range("C12").Value = 0

$ws.Range("B13").Value = 7
$ws.Range("C13").Value = 28

$ws.Range("B14").Value = 476
$ws.Range("C14").Value = 0

$ws.Range("B15").Value = 19
$ws.Range("C15").Value = 0
